$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C14").Value = 8502
$ws.Range("C15:C25").Value = 8477
$ws.Range("C26:C32").Value = 8339
$ws.Range("C33:C43").Value = 8112
$ws.Range("C44:C52").Value = 7800
$ws.Range("C53:C62").Value = 7750
$ws.Range("C63:C85").Value = 7748
$ws.Range("C86:C252").Value = 7293
